# Opp Test Data, Eng Detail, Add Counterparty - 10 Oct 2025
#
# Real content edit: the CAO user "Gemma Hardy" was swapped out for
# "Jennie Stewart" on the CAOUsers sheet.  Everything else in the target
# diff (shared-string re-indexing on every other sheet, the shrunk
# cellXfs table) is a mechanical side effect of that single value edit
# plus the user clearing stale formatting on one cell / switching the
# active sheet before saving, so we reproduce those user actions too.

$wb = $excel.ActiveWorkbook

# --- Content change ----------------------------------------------------
# CAOUsers!A2 : "Gemma Hardy" -> "Jennie Stewart"
$caoUsers = $wb.Worksheets.Item("CAOUsers")
$caoUsers.Range("A2").Value = "Jennie Stewart"

# --- Formatting cleanup --------------------------------------------------
# AddOpportunity!D2 carried an applied-but-redundant cell style (identical
# to the default). Clearing it drops the duplicate style definition.
$addOpportunity = $wb.Worksheets.Item("AddOpportunity")
$addOpportunity.Range("D2").ClearFormats() | Out-Null

# --- View / selection state ---------------------------------------------
# Selection left on AddOpportunity...
$addOpportunity.Range("C18").Select() | Out-Null

# ...then CAOUsers becomes the active/visible sheet with its own selection.
$caoUsers.Activate()
$caoUsers.Range("F18").Select() | Out-Null
